$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Status text: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (shared string reused by Overview!E2:F3 and the per-locale Status column)
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2. Per-locale handback details: zh-cn sheet
#    I = Latest Target File (hyperlink to the source .md)
#    J = Latest Handback File (generated xliff name)
#    K = Latest Handback DateTime
# ---------------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"

$zhcn.Range("I2").Value = "4e40cad4-ee76-4051-99d1-db7f748839f7.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d831c8246f15b59dc2ec087c931f0228c3529efb/e2e/4e40cad4-ee76-4051-99d1-db7f748839f7.md", "", "", "4e40cad4-ee76-4051-99d1-db7f748839f7.md") | Out-Null
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = 15570276
$zhcn.Range("J2").Value = "4e40cad4-ee76-4051-99d1-db7f748839f7.23306065a33af6ca392c7314c3e62436e6bfd177.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-19 19:10:36"

$zhcn.Range("I3").Value = "4e40cad4-ee76-4051-99d1-db7f748839f7.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d831c8246f15b59dc2ec087c931f0228c3529efb/e2e/4e40cad4-ee76-4051-99d1-db7f748839f7.md", "", "", "4e40cad4-ee76-4051-99d1-db7f748839f7.md") | Out-Null
$zhcn.Range("I3").Font.Underline = 2
$zhcn.Range("I3").Font.Color = 15570276
$zhcn.Range("J3").Value = "4e40cad4-ee76-4051-99d1-db7f748839f7.23306065a33af6ca392c7314c3e62436e6bfd177.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-19 19:10:36"

# widened columns on zh-cn
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# 3. Per-locale handback details: de-de sheet
# ---------------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("C3").Value = "Handed back: in sync with en-US"

$dede.Range("I2").Value = "4e40cad4-ee76-4051-99d1-db7f748839f7.md"
$dede.Hyperlinks.Add($dede.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d831c8246f15b59dc2ec087c931f0228c3529efb/e2e/4e40cad4-ee76-4051-99d1-db7f748839f7.md", "", "", "4e40cad4-ee76-4051-99d1-db7f748839f7.md") | Out-Null
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = 15570276
$dede.Range("J2").Value = "4e40cad4-ee76-4051-99d1-db7f748839f7.23306065a33af6ca392c7314c3e62436e6bfd177.de-de.xlf"
$dede.Range("K2").Value = "2016-08-19 19:10:43"

$dede.Range("I3").Value = "4e40cad4-ee76-4051-99d1-db7f748839f7.md"
$dede.Hyperlinks.Add($dede.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d831c8246f15b59dc2ec087c931f0228c3529efb/e2e/4e40cad4-ee76-4051-99d1-db7f748839f7.md", "", "", "4e40cad4-ee76-4051-99d1-db7f748839f7.md") | Out-Null
$dede.Range("I3").Font.Underline = 2
$dede.Range("I3").Font.Color = 15570276
$dede.Range("J3").Value = "4e40cad4-ee76-4051-99d1-db7f748839f7.23306065a33af6ca392c7314c3e62436e6bfd177.de-de.xlf"
$dede.Range("K3").Value = "2016-08-19 19:10:43"

# widened columns on de-de
$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# 4. Widened columns on Overview (E, F)
# ---------------------------------------------------------------------------
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527
